$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (it holds numeric-looking strings like
# "27.892.37" / "324.74") while writing new values, matching the original
# inlineStr cell type; ClearFormats() afterwards restores the default (no) style
# so only cell values change, not formatting.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.892.37"
$ws.Range("D3").Value = "1.871.33"
$ws.Range("D5").Value = "324.74"
$ws.Range("D7").Value = "0.4435"
$ws.Range("D8").Value = "0.3828"
$ws.Range("D9").Value = "0.07482"
$ws.Range("D10").Value = "0.8929"
$ws.Range("D12").Value = "1.884.11"
$ws.Range("D13").Value = "5.616"
$ws.Range("D14").Value = "6.792"
$ws.Range("D15").Value = "0.07219"
$ws.Range("D16").Value = "85.08"
$ws.Range("D17").Value = "1.038"
$ws.Range("D18").Value = "0.000009152"
$ws.Range("D20").Value = "15.64"
$ws.Range("D21").Value = "27.923.62"
$ws.Range("D22").Value = "5.344"
$ws.Range("D24").Value = "2.106.03"
$ws.Range("D25").Value = "2.025"
$ws.Range("D26").Value = "158.73"
$ws.Range("D27").Value = "18.94"
$ws.Range("D28").Value = "5.433"
$ws.Range("D29").Value = "1.998"
$ws.Range("D30").Value = "118.64"
$ws.Range("D31").Value = "0.09078"
$ws.Range("D32").Value = "1.241"
$ws.Range("D33").Value = "0.7845"
$ws.Range("D34").Value = "4.625"
$ws.Range("D35").Value = "3.019"
$ws.Range("D36").Value = "1.033"
$ws.Range("D37").Value = "1.147"
$ws.Range("D38").Value = "0.01996"
$ws.Range("D39").Value = "0.05383"
$ws.Range("D40").Value = "2.894"
$ws.Range("D41").Value = "0.5245"
$ws.Range("D42").Value = "0.1701"
$ws.Range("D43").Value = "6.924"
$ws.Range("D44").Value = "8.894"
$ws.Range("D45").Value = "112.32"
$ws.Range("D46").Value = "10.75"
$ws.Range("D47").Value = "0.06623"
$ws.Range("D48").Value = "1.034"
$ws.Range("D49").Value = "1.729"
$ws.Range("D50").Value = "0.4763"
$ws.Range("D51").Value = "1.924"

$priceRange.ClearFormats()

# Volume(1h) column values are percentage text (e.g. "  +1.98%  ") and are not
# auto-coerced to numbers by Excel, so plain Value assignment is sufficient.
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("E12").Value = "  -3.16%  "
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  +7.17%  "
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +3.81%  "
$ws.Range("E29").Value = "  +4.24%  "
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("E35").Value = "  +5.74%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("E44").Value = "  +5.41%  "
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("E47").Value = "  +5.32%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("E51").Value = "  +2.75%  "
